$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price and volume figures (cell values kept as text,
# matching the original inline-string cells: a leading apostrophe forces
# text entry, then the style is reset to "Normal" to avoid introducing a
# quote-prefix style that was not present in the source file).

$ws.Range("D2").Value = "'60.277.67"
$ws.Range("E2").Value = "'  +4.04%  "
$ws.Range("D3").Value = "'2.448.11"
$ws.Range("E4").Value = "'  +0.01%  "
$ws.Range("D5").Value = "'557.01"
$ws.Range("E5").Value = "'  +3.39%  "
$ws.Range("D6").Value = "'138.89"
$ws.Range("E6").Value = "'  +2.49%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "'  -0.04%  "
$ws.Range("E8").Value = "'  +1.39%  "
$ws.Range("D9").Value = "'0.108"
$ws.Range("D10").Value = "'5.82"
$ws.Range("E10").Value = "'  +5.04%  "
$ws.Range("E11").Value = "'  +2.37%  "
$ws.Range("E12").Value = "'  -1.89%  "
$ws.Range("D13").Value = "'25.00"
$ws.Range("E13").Value = "'  +5.08%  "
$ws.Range("D14").Value = "'2.883.55"
$ws.Range("E14").Value = "'  +4.25%  "
$ws.Range("D15").Value = "'60.122.24"
$ws.Range("E16").Value = "'  +5.47%  "
$ws.Range("D17").Value = "'2.463.64"
$ws.Range("E17").Value = "'  +5.17%  "
$ws.Range("D18").Value = "'11.48"
$ws.Range("E18").Value = "'  +7.51%  "
$ws.Range("E19").Value = "'  +3.99%  "
$ws.Range("D20").Value = "'335.86"
$ws.Range("E20").Value = "'  +1.25%  "
$ws.Range("D21").Value = "'6.92"
$ws.Range("E21").Value = "'  +2.25%  "
$ws.Range("E22").Value = "'  +0.07%  "
$ws.Range("E23").Value = "'  +3.17%  "
$ws.Range("D24").Value = "'0.171"
$ws.Range("E24").Value = "'  +2.52%  "
$ws.Range("D25").Value = "'8.60"
$ws.Range("E25").Value = "'  +1.21%  "
$ws.Range("E26").Value = "'  +0.10%  "
$ws.Range("E27").Value = "'  +0.70%  "
$ws.Range("D28").Value = "'0.0₃0799"
$ws.Range("E28").Value = "'  +8.48%  "
$ws.Range("E29").Value = "'  +3.88%  "
$ws.Range("E30").Value = "'  +3.16%  "
$ws.Range("D31").Value = "'171.18"
$ws.Range("E31").Value = "'  -0.44%  "
$ws.Range("D32").Value = "'18.86"
$ws.Range("E32").Value = "'  +2.21%  "
$ws.Range("E33").Value = "'  -0.13%  "
$ws.Range("E35").Value = "'  +6.22%  "
$ws.Range("D36").Value = "'4.31"
$ws.Range("E36").Value = "'  +2.11%  "
$ws.Range("E37").Value = "'  +0.09%  "
$ws.Range("E38").Value = "'  +0.51%  "
$ws.Range("D39").Value = "'40.14"
$ws.Range("E39").Value = "'  +2.22%  "
$ws.Range("D40").Value = "'0.419"
$ws.Range("E40").Value = "'  +10.96%  "
$ws.Range("D41").Value = "'316.31"
$ws.Range("E41").Value = "'  +6.74%  "
$ws.Range("E42").Value = "'  +2.43%  "
$ws.Range("D43").Value = "'144.09"
$ws.Range("E43").Value = "'  -1.15%  "
$ws.Range("E44").Value = "'  +1.72%  "
$ws.Range("D45").Value = "'0.0526"
$ws.Range("E45").Value = "'  +4.68%  "
$ws.Range("D46").Value = "'19.53"
$ws.Range("E46").Value = "'  +1.40%  "
$ws.Range("D47").Value = "'0.412"
$ws.Range("E47").Value = "'  +6.90%  "
$ws.Range("D48").Value = "'0.576"
$ws.Range("E48").Value = "'  +2.36%  "
$ws.Range("E49").Value = "'  +2.93%  "
$ws.Range("E50").Value = "'  -0.28%  "
$ws.Range("E51").Value = "'  +5.24%  "

# Reset style to the workbook default so no quotePrefix / number-format
# style gets attached to these cells (they had no explicit style before).
$cells = @("D2","E2","D3","E4","D5","E5","D6","E6","D7","E7","E8","D9","D10","E10","E11","E12","D13","E13","D14","E14","D15","E16","D17","E17","D18","E18","E19","D20","E20","D21","E21","E22","E23","D24","E24","D25","E25","E26","E27","D28","E28","E29","E30","D31","E31","D32","E32","E33","E35","D36","E36","E37","E38","D39","E39","D40","E40","D41","E41","E42","D43","E43","E44","D45","E45","D46","E46","D47","E47","D48","E48","E49","E50","E51")
foreach ($cellRef in $cells) {
    $ws.Range($cellRef).Style = "Normal"
}

Write-Host "Applied cryptos update"
